$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.4859026666666666
$ws.Range("H2").Value = 1.457708
$ws.Range("M2").Value = 0.106124
$ws.Range("N2").Value = 0.318372
$ws.Range("O2").Value = 0.08094716512538251
$ws.Range("P2").Value = 0.08094716512538253
$ws.Range("Q2").Value = 0.05156593459733333
$ws.Range("R2").Value = 0.464093411376
$ws.Range("S2").Value = 0.08094716512538251
$ws.Range("T2").Value = 0.08094716512538253

# Row 3
$ws.Range("G3").Value = 0.4859026666666666
$ws.Range("H3").Value = 1.457708
$ws.Range("O3").Value = 0.8331551016962769
$ws.Range("P3").Value = 0.833155101696277
$ws.Range("Q3").Value = 0.5307464618057778
$ws.Range("R3").Value = 4.776718156252
$ws.Range("S3").Value = 0.8331551016962769
$ws.Range("T3").Value = 0.833155101696277

# Row 4
$ws.Range("G4").Value = 0.4859026666666666
$ws.Range("H4").Value = 1.457708
$ws.Range("M4").Value = 0.1126143333333333
$ws.Range("N4").Value = 0.337843
$ws.Range("O4").Value = 0.08589773317834044
$ws.Range("P4").Value = 0.08589773317834046
$ws.Range("Q4").Value = 0.05471960487155555
$ws.Range("R4").Value = 0.492476443844
$ws.Range("S4").Value = 0.08589773317834044
$ws.Range("T4").Value = 0.08589773317834046
